$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.554.11"
$ws.Range("E2").Value = "  +0.51%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.841.47"
$ws.Range("E3").Value = "  +0.02%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.07%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'259.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.03%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  +0.07%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +0.80%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.3181"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.44%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "'0.06795"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.36%  "

# Row 10 - Solana
$ws.Range("D10").Value = "'18.76"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.06%  "

# Row 11 - Polygon
$ws.Range("D11").Value = "'0.7811"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.18%  "

# Row 12 - TRON
$ws.Range("D12").Value = "'0.07770"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.90%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.821.08"
$ws.Range("E13").Value = "  -1.06%  "

# Row 14 - Litecoin
$ws.Range("D14").Value = "'87.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.13%  "

# Row 15 - Polkadot
$ws.Range("D15").Value = "'5.016"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.11%  "

# Row 16 - BinanceUSD
$ws.Range("D16").Value = "'1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.00%  "

# Row 17 - Avalanche
$ws.Range("E17").Value = "  -0.03%  "

# Row 18 - Dai
$ws.Range("E18").Value = "  -0.06%  "

# Row 19 - ShibaInu
$ws.Range("D19").Value = "'0.000007942"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.10%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "26.575.18"
$ws.Range("E20").Value = "  +0.48%  "

# Row 21 - WrappedliquidstakedEther2.0
$ws.Range("D21").Value = "2.060.94"
$ws.Range("E21").Value = "  -0.63%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'4.610"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.17%  "

# Row 23 - Chainlink
$ws.Range("D23").Value = "'5.973"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.42%  "

# Row 24 - Cosmos
$ws.Range("D24").Value = "'9.346"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.21%  "

# Row 25 - LidoDAOToken (was Monero)
$ws.Range("B25").Value = "LidoDAOToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"

$ws.Range("D25").Value = "'2.224"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.49%  "

# Row 26 - Monero (was LidoDAOToken)
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"

$ws.Range("D26").Value = "'142.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.22%  "

# Row 27 - Toncoin
$ws.Range("D27").Value = "'1.673"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.89%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "'16.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.41%  "

# Row 29 - BitcoinCash
$ws.Range("D29").Value = "'111.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.76%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Range("D30").Value = "'4.183"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.62%  "

# Row 31 - Stellar
$ws.Range("D31").Value = "'0.08731"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.26%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "'4.076"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.02%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "'0.04888"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.36%  "

# Row 34 - ImmutableX
$ws.Range("D34").Value = "'0.7237"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.21%  "

# Row 35 - ARBITRUM
$ws.Range("E35").Value = "  +0.82%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  +0.30%  "

# Row 37 - MXToken
$ws.Range("D37").Value = "'3.095"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.38%  "

# Row 38 - RenderToken
$ws.Range("D38").Value = "'2.231"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.22%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  -0.24%  "

# Row 40 - TheSandbox
$ws.Range("D40").Value = "'0.4823"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.27%  "

# Row 41 - TrustWalletToken
$ws.Range("D41").Value = "'0.8985"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.65%  "

# Row 42 - Quant
$ws.Range("D42").Value = "'110.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.73%  "

# Row 43 - FraxShare
$ws.Range("D43").Value = "'5.921"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.34%  "

# Row 44 - PaxDollar
$ws.Range("E44").Value = "  +0.12%  "

# Row 45 - Aptos
$ws.Range("D45").Value = "'7.649"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.05%  "

# Row 46 - Decentraland
$ws.Range("D46").Value = "'0.4171"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.64%  "

# Row 47 - EnergySwap
$ws.Range("D47").Value = "'9.007"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.59%  "

# Row 48 - Algorand
$ws.Range("E48").Value = "  +1.42%  "

# Row 49 - Cronos
$ws.Range("E49").Value = "  -0.52%  "

# Row 50 - Elrond
$ws.Range("D50").Value = "'34.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.29%  "

# Row 51 - EOS
$ws.Range("D51").Value = "'0.8935"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.84%  "

